# Tutorial 6: update attendance date strings from DD/MM/YYYY to DD-MM-YYYY
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Range("A$row")
    # Force text interpretation so day/month-ambiguous strings (e.g. "01-08-2022")
    # are not auto-converted into a date serial by the COM value setter.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    # Drop the temporary Text number format again so the cell's style matches
    # the original (unstyled) cell - only the literal text should change.
    $cell.ClearFormats()
}
